$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.372.41'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.645.24'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').Value = '2.645.91'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = '3.131.68'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000186'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').Value = '68.387.35'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '2.674.65'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '363.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('E23').Value = '  -3.22%  '
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').Value = '2.778.23'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '561.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.73'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('E44').Value = '  -4.47%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0773'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('E51').Value = '  +0.91%  '
